# Preferential queue and preferential queue tests
#
# Adds a new "205" time-bucket (row 42) to every results sheet, and
# refreshes the "network" counter series (and its mirrored copy in the
# D column of the per-distributor sheets) with the updated values that
# come along with that extra data point.

$wb = $excel.ActiveWorkbook

# New cumulative "network" values for rows 2..42 (time = 5..205).
# These replace the previous B/C values on the "network" sheet, and the
# previous D values on the "0TestDistributor" / "1TestDistributor" sheets.
$networkValues = @(17, 28, 38, 51, 64, 79, 98, 113, 121, 131, 141, 146, 155, 170, 180, 185, 191, 199, 206, 214, 224, 236, 249, 257, 261, 264, 275, 289, 300, 310, 318, 325, 335, 346, 355, 363, 371, 379, 388, 395, 396)

# ---- success ----------------------------------------------------------
$ws = $wb.Worksheets.Item("success")
$ws.Cells.Item(42, 1).Value = 205
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 0

# ---- fail ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("fail")
$ws.Cells.Item(42, 1).Value = 205
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 0

# ---- network --------------------------------------------------------
$ws = $wb.Worksheets.Item("network")
for ($i = 0; $i -lt $networkValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $networkValues[$i]
    $ws.Cells.Item($row, 3).Value = $networkValues[$i]
}
$ws.Cells.Item(42, 1).Value = 205

# ---- 0TestDistributor -------------------------------------------------
$ws = $wb.Worksheets.Item("0TestDistributor")
for ($i = 0; $i -lt $networkValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $networkValues[$i]
}
$ws.Cells.Item(42, 1).Value = 205
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 0

# ---- 1TestDistributor -------------------------------------------------
$ws = $wb.Worksheets.Item("1TestDistributor")
for ($i = 0; $i -lt $networkValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $networkValues[$i]
}
$ws.Cells.Item(42, 1).Value = 205
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 0
